$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.PageSetup.PrintArea = "$A$1:$BB$15"
